# Applies the numeric corrections to the leve-profit tables (columns H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR worksheets, as
# produced by the scheduled market-price refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 3783.75
$ws.Range("J6").Value = 9999
$ws.Range("L6").Value = 29997
$ws.Range("N6").Value = -30221
# Row 100
$ws.Range("H100").Value = 1528.1
$ws.Range("J100").Value = 2102
$ws.Range("L100").Value = 2102
$ws.Range("N100").Value = -3184
# Row 132
$ws.Range("H132").Value = 2013.8572
$ws.Range("I132").Value = 1518.2
$ws.Range("J132").Value = 4987.8
$ws.Range("K132").Value = 4554.6
$ws.Range("L132").Value = 14963.4
$ws.Range("M132").Value = -2024.6
$ws.Range("N132").Value = -20023.4

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3280.0454
$ws.Range("I2").Value = 2733.05
$ws.Range("J2").Value = 8750
$ws.Range("K2").Value = 2733.05
$ws.Range("L2").Value = 8750
$ws.Range("M2").Value = -2620.05
$ws.Range("N2").Value = -8976
# Row 22
$ws.Range("H22").Value = 8553.714
# Row 41
$ws.Range("H41").Value = 21369
$ws.Range("I41").Value = 922.5
$ws.Range("J41").Value = 35000
$ws.Range("K41").Value = 922.5
$ws.Range("L41").Value = 35000
$ws.Range("M41").Value = -508.5
$ws.Range("N41").Value = -35828
# Row 116
$ws.Range("H116").Value = 3280.0454
$ws.Range("I116").Value = 2733.05
$ws.Range("J116").Value = 8750
$ws.Range("K116").Value = 2733.05
$ws.Range("L116").Value = 8750
$ws.Range("M116").Value = -439.0500000000002
$ws.Range("N116").Value = -13338
# Row 122
$ws.Range("H122").Value = 1517.875
$ws.Range("I122").Value = 1517.875
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4553.625
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2103.625
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3280.0454
$ws.Range("I3").Value = 2733.05
$ws.Range("J3").Value = 8750
$ws.Range("K3").Value = 2733.05
$ws.Range("L3").Value = 8750
$ws.Range("M3").Value = -2619.05
$ws.Range("N3").Value = -8978
# Row 82
$ws.Range("H82").Value = 20700
$ws.Range("I82").Value = 6000
$ws.Range("K82").Value = 6000
$ws.Range("M82").Value = -5617
# Row 85
$ws.Range("H85").Value = 20700
$ws.Range("I85").Value = 6000
$ws.Range("K85").Value = 6000
$ws.Range("M85").Value = -4674
# Row 86
$ws.Range("H86").Value = 1926.5385
$ws.Range("I86").Value = 1839.2222
$ws.Range("J86").Value = 2123
$ws.Range("K86").Value = 1839.2222
$ws.Range("L86").Value = 2123
$ws.Range("M86").Value = -716.2221999999999
$ws.Range("N86").Value = -4369
# Row 89
$ws.Range("H89").Value = 1926.5385
$ws.Range("I89").Value = 1839.2222
$ws.Range("J89").Value = 2123
$ws.Range("K89").Value = 9196.110999999999
$ws.Range("L89").Value = 10615
$ws.Range("M89").Value = -3580.110999999999
$ws.Range("N89").Value = -21847

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 90917010
$ws.Range("J31").Value = 10910
$ws.Range("L31").Value = 10910
$ws.Range("N31").Value = -11500
# Row 34
$ws.Range("H34").Value = 90917010
$ws.Range("J34").Value = 10910
$ws.Range("L34").Value = 10910
$ws.Range("N34").Value = -11314
# Row 41
$ws.Range("H41").Value = 20591.666
$ws.Range("J41").Value = 20591.666
$ws.Range("L41").Value = 20591.666
$ws.Range("N41").Value = -21447.666
# Row 58
$ws.Range("H58").Value = 8695.044
$ws.Range("J58").Value = 11605.286
$ws.Range("L58").Value = 11605.286
$ws.Range("N58").Value = -12011.286
# Row 59
$ws.Range("H59").Value = 14000
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
# Row 136
$ws.Range("H136").Value = 8695.044
$ws.Range("J136").Value = 11605.286
$ws.Range("L136").Value = 34815.858
$ws.Range("N136").Value = -39915.858

$ws = $wb.Worksheets.Item("CUL")
# Row 29
$ws.Range("H29").Value = 75433.336
$ws.Range("I29").Value = 650
$ws.Range("J29").Value = 225000
$ws.Range("K29").Value = 1950
$ws.Range("L29").Value = 675000
$ws.Range("M29").Value = -1673
$ws.Range("N29").Value = -675554
# Row 39
$ws.Range("H39").Value = 1804.875
$ws.Range("I39").Value = 1479.6666
$ws.Range("K39").Value = 4438.9998
$ws.Range("M39").Value = -4144.9998
# Row 50
$ws.Range("H50").Value = 717.9231
$ws.Range("I50").Value = 2105
$ws.Range("K50").Value = 6315
$ws.Range("M50").Value = -5834
# Row 53
$ws.Range("H53").Value = 717.9231
$ws.Range("I53").Value = 2105
$ws.Range("K53").Value = 6315
$ws.Range("M53").Value = -5834
# Row 55
$ws.Range("H55").Value = 94621.875
$ws.Range("J55").Value = 150798
$ws.Range("L55").Value = 452394
$ws.Range("N55").Value = -452748
# Row 104
$ws.Range("H104").Value = 7999
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()
# Row 137
$ws.Range("H137").Value = 10736
$ws.Range("J137").Value = 24500
$ws.Range("L137").Value = 73500
$ws.Range("N137").Value = -83700

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 9113.541999999999
$ws.Range("I70").Value = 6902.353
$ws.Range("K70").Value = 6902.353
$ws.Range("M70").Value = -6632.353
# Row 73
$ws.Range("H73").Value = 9113.541999999999
$ws.Range("I73").Value = 6902.353
$ws.Range("K73").Value = 6902.353
$ws.Range("M73").Value = -5966.353
# Row 136
$ws.Range("H136").Value = 34480.684
$ws.Range("J136").Value = 34480.684
$ws.Range("L136").Value = 103442.052
$ws.Range("N136").Value = -108542.052

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 7103.4644
$ws.Range("J46").Value = 6755.4443
$ws.Range("L46").Value = 6755.4443
$ws.Range("N46").Value = -7131.4443
# Row 136
$ws.Range("H136").Value = 6622.136
$ws.Range("I136").Value = 6619.512
$ws.Range("K136").Value = 19858.536
$ws.Range("M136").Value = -17308.536

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 2999.8333
$ws.Range("I2").Value = 2999
$ws.Range("K2").Value = 2999
$ws.Range("M2").Value = -2887
# Row 54
$ws.Range("H54").Value = 99533.336
$ws.Range("J54").Value = 99533.336
$ws.Range("L54").Value = 99533.336
$ws.Range("N54").Value = -100573.336
# Row 96
$ws.Range("H96").Value = 1201.25
$ws.Range("I96").Value = 1103
$ws.Range("J96").Value = 1299.5
$ws.Range("K96").Value = 1103
$ws.Range("L96").Value = 1299.5
$ws.Range("M96").Value = 270
$ws.Range("N96").Value = -4045.5
# Row 136
$ws.Range("H136").Value = 2795.7368
$ws.Range("I136").Value = 2117.7222
$ws.Range("K136").Value = 6353.1666
$ws.Range("M136").Value = -3803.1666
